$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object "object[,]" 1,20
$row2[0,0] = "ECs"
$row2[0,1] = "Pdgfb"
$row2[0,2] = "Art1"
$row2[0,3] = "ECs"
$row2[0,4] = 3
$row2[0,5] = 1
$row2[0,6] = 29.253501
$row2[0,7] = 87.760503
$row2[0,8] = 0.7876335333413836
$row2[0,9] = 0.7876335333413838
$row2[0,10] = 1
$row2[0,11] = 0.3333333333333333
$row2[0,12] = 0.04418433333333333
$row2[0,13] = 0.132553
$row2[0,14] = 0.01996731490055206
$row2[0,15] = 0.01996731490055207
$row2[0,16] = 1.292546439351
$row2[0,17] = 11.632917954159
$row2[0,18] = 0.01572692678646188
$row2[0,19] = 0.01572692678646188
$ws.Range("A2:T2").Value = $row2

$row3 = New-Object "object[,]" 1,20
$row3[0,0] = "ECs"
$row3[0,1] = "Pdgfb"
$row3[0,2] = "Art1"
$row3[0,3] = "FAPs"
$row3[0,4] = 3
$row3[0,5] = 1
$row3[0,6] = 29.253501
$row3[0,7] = 87.760503
$row3[0,8] = 0.7876335333413836
$row3[0,9] = 0.7876335333413838
$row3[0,10] = 3
$row3[0,11] = 1
$row3[0,12] = 0.5088386666666667
$row3[0,13] = 1.526516
$row3[0,14] = 0.2299489688858882
$row3[0,15] = 0.2299489688858882
$row3[0,16] = 14.885312444172
$row3[0,17] = 133.967811997548
$row3[0,18] = 0.1811155188518
$row3[0,19] = 0.1811155188518001
$ws.Range("A3:T3").Value = $row3

$row4 = New-Object "object[,]" 1,20
$row4[0,0] = "ECs"
$row4[0,1] = "Pdgfb"
$row4[0,2] = "Art1"
$row4[0,3] = "M2"
$row4[0,4] = 3
$row4[0,5] = 1
$row4[0,6] = 29.253501
$row4[0,7] = 87.760503
$row4[0,8] = 0.7876335333413836
$row4[0,9] = 0.7876335333413838
$row4[0,10] = 1
$row4[0,11] = 0.3333333333333333
$row4[0,12] = 0.007129666666666666
$row4[0,13] = 0.021389
$row4[0,14] = 0.003221963278144652
$row4[0,15] = 0.003221963278144652
$row4[0,16] = 0.208567710963
$row4[0,17] = 1.877109398667
$row4[0,18] = 0.00253772632106126
$row4[0,19] = 0.00253772632106126
$ws.Range("A4:T4").Value = $row4

$row5 = New-Object "object[,]" 1,20
$row5[0,0] = "ECs"
$row5[0,1] = "Pdgfb"
$row5[0,2] = "Art1"
$row5[0,3] = "sCs"
$row5[0,4] = 3
$row5[0,5] = 1
$row5[0,6] = 29.253501
$row5[0,7] = 87.760503
$row5[0,8] = 0.7876335333413836
$row5[0,9] = 0.7876335333413838
$row5[0,10] = 3
$row5[0,11] = 1
$row5[0,12] = 1.652680333333334
$row5[0,13] = 4.958041000000001
$row5[0,14] = 0.746861752935415
$row5[0,15] = 0.7468617529354151
$row5[0,16] = 48.34668578384701
$row5[0,17] = 435.120172054623
$row5[0,18] = 0.5882533613820604
$row5[0,19] = 0.5882533613820606
$ws.Range("A5:T5").Value = $row5

$row6 = New-Object "object[,]" 1,20
$row6[0,0] = "M2"
$row6[0,1] = "Pdgfb"
$row6[0,2] = "Art1"
$row6[0,3] = "ECs"
$row6[0,4] = 3
$row6[0,5] = 1
$row6[0,6] = 4.749137666666667
$row6[0,7] = 14.247413
$row6[0,8] = 0.1278677748937237
$row6[0,9] = 0.1278677748937237
$row6[0,10] = 1
$row6[0,11] = 0.3333333333333333
$row6[0,12] = 0.04418433333333333
$row6[0,13] = 0.132553
$row6[0,14] = 0.01996731490055206
$row6[0,15] = 0.01996731490055207
$row6[0,16] = 0.2098374817098889
$row6[0,17] = 1.888537335389
$row6[0,18] = 0.002553176126935886
$row6[0,19] = 0.002553176126935887
$ws.Range("A6:T6").Value = $row6

$row7 = New-Object "object[,]" 1,20
$row7[0,0] = "M2"
$row7[0,1] = "Pdgfb"
$row7[0,2] = "Art1"
$row7[0,3] = "FAPs"
$row7[0,4] = 3
$row7[0,5] = 1
$row7[0,6] = 4.749137666666667
$row7[0,7] = 14.247413
$row7[0,8] = 0.1278677748937237
$row7[0,9] = 0.1278677748937237
$row7[0,10] = 3
$row7[0,11] = 1
$row7[0,12] = 0.5088386666666667
$row7[0,13] = 1.526516
$row7[0,14] = 0.2299489688858882
$row7[0,15] = 0.2299489688858882
$row7[0,16] = 2.416544878123111
$row7[0,17] = 21.748903903108
$row7[0,18] = 0.02940306299054462
$row7[0,19] = 0.02940306299054463
$ws.Range("A7:T7").Value = $row7

$row8 = New-Object "object[,]" 1,20
$row8[0,0] = "M2"
$row8[0,1] = "Pdgfb"
$row8[0,2] = "Art1"
$row8[0,3] = "M2"
$row8[0,4] = 3
$row8[0,5] = 1
$row8[0,6] = 4.749137666666667
$row8[0,7] = 14.247413
$row8[0,8] = 0.1278677748937237
$row8[0,9] = 0.1278677748937237
$row8[0,10] = 1
$row8[0,11] = 0.3333333333333333
$row8[0,12] = 0.007129666666666666
$row8[0,13] = 0.021389
$row8[0,14] = 0.003221963278144652
$row8[0,15] = 0.003221963278144652
$row8[0,16] = 0.03385976851744445
$row8[0,17] = 0.304737916657
$row8[0,18] = 0.0004119852751656445
$row8[0,19] = 0.0004119852751656446
$ws.Range("A8:T8").Value = $row8

$row9 = New-Object "object[,]" 1,20
$row9[0,0] = "M2"
$row9[0,1] = "Pdgfb"
$row9[0,2] = "Art1"
$row9[0,3] = "sCs"
$row9[0,4] = 3
$row9[0,5] = 1
$row9[0,6] = 4.749137666666667
$row9[0,7] = 14.247413
$row9[0,8] = 0.1278677748937237
$row9[0,9] = 0.1278677748937237
$row9[0,10] = 3
$row9[0,11] = 1
$row9[0,12] = 1.652680333333334
$row9[0,13] = 4.958041000000001
$row9[0,14] = 0.746861752935415
$row9[0,15] = 0.7468617529354151
$row9[0,16] = 7.848806421992557
$row9[0,17] = 70.63925779793301
$row9[0,18] = 0.09549955050107753
$row9[0,19] = 0.09549955050107757
$ws.Range("A9:T9").Value = $row9

$row10 = New-Object "object[,]" 1,20
$row10[0,0] = "sCs"
$row10[0,1] = "Pdgfb"
$row10[0,2] = "Art1"
$row10[0,3] = "ECs"
$row10[0,4] = 3
$row10[0,5] = 1
$row10[0,6] = 3.138366333333334
$row10[0,7] = 9.415099000000001
$row10[0,8] = 0.08449869176489255
$row10[0,9] = 0.08449869176489258
$row10[0,10] = 1
$row10[0,11] = 0.3333333333333333
$row10[0,12] = 0.04418433333333333
$row10[0,13] = 0.132553
$row10[0,14] = 0.01996731490055206
$row10[0,15] = 0.01996731490055207
$row10[0,16] = 0.1386666241941111
$row10[0,17] = 1.247999617747
$row10[0,18] = 0.001687211987154295
$row10[0,19] = 0.001687211987154296
$ws.Range("A10:T10").Value = $row10

$row11 = New-Object "object[,]" 1,20
$row11[0,0] = "sCs"
$row11[0,1] = "Pdgfb"
$row11[0,2] = "Art1"
$row11[0,3] = "FAPs"
$row11[0,4] = 3
$row11[0,5] = 1
$row11[0,6] = 3.138366333333334
$row11[0,7] = 9.415099000000001
$row11[0,8] = 0.08449869176489255
$row11[0,9] = 0.08449869176489258
$row11[0,10] = 3
$row11[0,11] = 1
$row11[0,12] = 0.5088386666666667
$row11[0,13] = 1.526516
$row11[0,14] = 0.2299489688858882
$row11[0,15] = 0.2299489688858882
$row11[0,16] = 1.596922140564889
$row11[0,17] = 14.372299265084
$row11[0,18] = 0.01943038704354353
$row11[0,19] = 0.01943038704354354
$ws.Range("A11:T11").Value = $row11

$row12 = New-Object "object[,]" 1,20
$row12[0,0] = "sCs"
$row12[0,1] = "Pdgfb"
$row12[0,2] = "Art1"
$row12[0,3] = "M2"
$row12[0,4] = 3
$row12[0,5] = 1
$row12[0,6] = 3.138366333333334
$row12[0,7] = 9.415099000000001
$row12[0,8] = 0.08449869176489255
$row12[0,9] = 0.08449869176489258
$row12[0,10] = 1
$row12[0,11] = 0.3333333333333333
$row12[0,12] = 0.007129666666666666
$row12[0,13] = 0.021389
$row12[0,14] = 0.003221963278144652
$row12[0,15] = 0.003221963278144652
$row12[0,16] = 0.02237550583455556
$row12[0,17] = 0.201379552511
$row12[0,18] = 0.0002722516819177477
$row12[0,19] = 0.0002722516819177478
$ws.Range("A12:T12").Value = $row12

$row13 = New-Object "object[,]" 1,20
$row13[0,0] = "sCs"
$row13[0,1] = "Pdgfb"
$row13[0,2] = "Art1"
$row13[0,3] = "sCs"
$row13[0,4] = 3
$row13[0,5] = 1
$row13[0,6] = 3.138366333333334
$row13[0,7] = 9.415099000000001
$row13[0,8] = 0.08449869176489255
$row13[0,9] = 0.08449869176489258
$row13[0,10] = 3
$row13[0,11] = 1
$row13[0,12] = 1.652680333333334
$row13[0,13] = 4.958041000000001
$row13[0,14] = 0.746861752935415
$row13[0,15] = 0.7468617529354151
$row13[0,16] = 5.186716317895446
$row13[0,17] = 46.68044686105901
$row13[0,18] = 0.06310884105227697
$row13[0,19] = 0.063108841052277
$ws.Range("A13:T13").Value = $row13

Write-Output "done"